$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextCell "D2" '40.105.04'
Set-TextCell "E2" '  +1.83%  '

Set-TextCell "D3" '2.203.43'
Set-TextCell "E3" '  +2.12%  '

Set-TextCell "E4" '  -0.16%  '

Set-TextCell "D5" '228.65'
Set-TextCell "E5" '  -0.15%  '

Set-TextCell "D6" '0.630'
Set-TextCell "E6" '  +1.35%  '

Set-TextCell "D7" '63.82'
Set-TextCell "E7" '  +1.25%  '

Set-TextCell "E8" '  +0.00%  '

Set-TextCell "D9" '0.398'
Set-TextCell "E9" '  +0.57%  '

Set-TextCell "D10" '0.0861'
Set-TextCell "E10" '  -0.59%  '

Set-TextCell "E11" '  +0.28%  '

Set-TextCell "D12" '16.07'
Set-TextCell "E12" '  +0.54%  '

Set-TextCell "D13" '2.530.98'
Set-TextCell "E13" '  +2.08%  '

Set-TextCell "D14" '22.19'
Set-TextCell "E14" '  +0.29%  '

Set-TextCell "D15" '0.821'
Set-TextCell "E15" '  +0.74%  '

Set-TextCell "D16" '5.59'
Set-TextCell "E16" '  +0.40%  '

Set-TextCell "D17" '2.204.09'
Set-TextCell "E17" '  +2.19%  '

Set-TextCell "D18" '40.014.12'
Set-TextCell "E18" '  +1.71%  '

Set-TextCell "D19" '0.0₃0910'
Set-TextCell "E19" '  +6.79%  '

Set-TextCell "D20" '72.39'
Set-TextCell "E20" '  +0.21%  '

Set-TextCell "D21" '6.08'
Set-TextCell "E21" '  -0.69%  '

Set-TextCell "D22" '231.92'
Set-TextCell "E22" '  +1.51%  '

Set-TextCell "E23" '  -0.01%  '

Set-TextCell "E24" '  +2.54%  '

Set-TextCell "D25" '2.37'
Set-TextCell "E25" '  +0.46%  '

Set-TextCell "D26" '9.68'
Set-TextCell "E26" '  -0.67%  '

Set-TextCell "D27" '171.82'
Set-TextCell "E27" '  -0.02%  '

Set-TextCell "E28" '  +1.72%  '

Set-TextCell "E29" '  +3.27%  '

Set-TextCell "D30" '20.14'
Set-TextCell "E30" '  +2.38%  '

Set-TextCell "E31" '  +6.01%  '

Set-TextCell "E32" '  +1.22%  '

Set-TextCell "D33" '4.60'
Set-TextCell "E33" '  -1.19%  '

Set-TextCell "D34" '4.74'
Set-TextCell "E34" '  -1.09%  '

Set-TextCell "D35" '7.07'
Set-TextCell "E35" '  -0.15%  '

Set-TextCell "D36" '0.0625'
Set-TextCell "E36" '  +0.69%  '

Set-TextCell "E37" '  +8.43%  '

Set-TextCell "D38" '2.46'
Set-TextCell "E38" '  +1.29%  '

Set-TextCell "D39" '0.999'
Set-TextCell "E39" '  -0.23%  '

Set-TextCell "D40" '4.99'
Set-TextCell "E40" '  +18.38%  '

Set-TextCell "D41" '103.59'
Set-TextCell "E41" '  -0.14%  '

Set-TextCell "D42" '0.0230'
Set-TextCell "E42" '  -0.41%  '

Set-TextCell "D43" '17.80'
Set-TextCell "E43" '  -1.35%  '

Set-TextCell "B44" 'TrustWalletToken'
Set-TextCell "C44" 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell "D44" '1.23'
Set-TextCell "E44" '  +3.14%  '

Set-TextCell "B45" 'Maker'
Set-TextCell "C45" 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell "D45" '1.519.51'
Set-TextCell "E45" '  -0.95%  '

Set-TextCell "D46" '8.28'
Set-TextCell "E46" '  +6.53%  '

Set-TextCell "E47" '  +0.55%  '

Set-TextCell "D48" '0.0927'
Set-TextCell "E48" '  -0.46%  '

Set-TextCell "D49" '2.80'
Set-TextCell "E49" '  -0.31%  '

Set-TextCell "D50" '0.000196'
Set-TextCell "E50" '  +33.47%  '

Set-TextCell "D51" '2.409.30'
Set-TextCell "E51" '  +1.93%  '
